# "Add files via upload" — workbook re-saved from Excel.
# Functional change: D1 no longer holds the volatile LOOKUP formula,
# just the literal phone number it had last resolved to. The active
# selection also moved from D12 back to D1.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

# Replace the formula in D1 with its last computed literal value.
$ws.Range("D1").Value = 5542999203443

# Restore the saved cursor/selection position to D1.
$ws.Range("D1").Select()
